# Auto commit at 2025-12-11  7:42:06.91
# Update the "Metrics" source values (propagates via formulas into the
# "today" sheet), then restore the original selections/active sheet.

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")
$metrics.Range("B2").Value  = 133267.69
$metrics.Range("B3").Value  = 114513.61000000002
$metrics.Range("B4").Value  = 40956.639999999999
$metrics.Range("B5").Value  = 5437
$metrics.Range("B6").Value  = 5335974.8000000007
$metrics.Range("B7").Value  = 4514866.57
$metrics.Range("B8").Value  = 1572913.5200000005
$metrics.Range("B9").Value  = 208144
$metrics.Range("B10").Value = 33801355.789999992
$metrics.Range("B11").Value = 31790141.73
$metrics.Range("B12").Value = 11854635.559999995
$metrics.Range("B13").Value = 1305774

# Move the selection on Metrics to D9->E15 as captured in the workbook view.
$metrics.Range("E15").Select() | Out-Null

# The "today" sheet is the tab that is active/selected in the saved workbook;
# select it last (and move its cursor to L22) so it ends up the active sheet.
$today = $wb.Worksheets.Item("today")
$today.Range("L22").Select() | Out-Null
